$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 166666670
$ws.Range("I100").Value = 166666670
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 166666670
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -166666129
$ws.Range("N100").ClearContents()
$ws.Range("H107").Value = 62513504
$ws.Range("I107").Value = 83334670
$ws.Range("J107").Value = 50000
$ws.Range("K107").Value = 83334670
$ws.Range("L107").Value = 50000
$ws.Range("M107").Value = -83332750
$ws.Range("N107").Value = -53840
$ws.Range("H129").Value = 791.93024
$ws.Range("I129").Value = 386.2
$ws.Range("J129").Value = 914.8788
$ws.Range("K129").Value = 1158.6
$ws.Range("L129").Value = 2744.6364
$ws.Range("M129").Value = 3841.4
$ws.Range("N129").Value = -12744.6364
$ws.Range("H132").Value = 2142.7585
$ws.Range("I132").Value = 2072.25
$ws.Range("J132").Value = 2481.2
$ws.Range("K132").Value = 6216.75
$ws.Range("L132").Value = 7443.599999999999
$ws.Range("M132").Value = -3686.75
$ws.Range("N132").Value = -12503.6
$ws.Range("H137").Value = 1908.0476
$ws.Range("I137").Value = 1935.5454
$ws.Range("J137").Value = 1877.8
$ws.Range("K137").Value = 5806.6362
$ws.Range("L137").Value = 5633.4
$ws.Range("M137").Value = -3256.6362
$ws.Range("N137").Value = -10733.4
$ws.Range("H138").Value = 19921.46
$ws.Range("I138").Value = 1060.5555
$ws.Range("J138").Value = 27465.822
$ws.Range("K138").Value = 3181.6665
$ws.Range("L138").Value = 82397.466
$ws.Range("M138").Value = 1958.3335
$ws.Range("N138").Value = -92677.466

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2355.3125
$ws.Range("I2").Value = 2400
$ws.Range("J2").Value = 2221.25
$ws.Range("K2").Value = 2400
$ws.Range("L2").Value = 2221.25
$ws.Range("M2").Value = -2287
$ws.Range("N2").Value = -2447.25
$ws.Range("H32").Value = 6600.5957
$ws.Range("I32").Value = 6249.032
$ws.Range("J32").Value = 7281.75
$ws.Range("K32").Value = 6249.032
$ws.Range("L32").Value = 7281.75
$ws.Range("M32").Value = -5962.032
$ws.Range("N32").Value = -7855.75
$ws.Range("H102").Value = 18519812
$ws.Range("I102").Value = 18519812
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 18519812
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -18518190
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value = 871.1
$ws.Range("I110").Value = 890.1111
$ws.Range("J110").Value = 700
$ws.Range("K110").Value = 890.1111
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = 1154.8889
$ws.Range("N110").Value = -4790
$ws.Range("H116").Value = 2355.3125
$ws.Range("I116").Value = 2400
$ws.Range("J116").Value = 2221.25
$ws.Range("K116").Value = 2400
$ws.Range("L116").Value = 2221.25
$ws.Range("M116").Value = -106
$ws.Range("N116").Value = -6809.25
$ws.Range("H122").Value = 2138212.5
$ws.Range("I122").Value = 2332413.5
$ws.Range("K122").Value = 6997240.5
$ws.Range("M122").Value = -6994790.5
$ws.Range("H132").Value = 4177.8423
$ws.Range("I132").Value = 1537.7368
$ws.Range("J132").Value = 9458.053
$ws.Range("K132").Value = 4613.2104
$ws.Range("L132").Value = 28374.159
$ws.Range("M132").Value = -2083.2104
$ws.Range("N132").Value = -33434.159
$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2355.3125
$ws.Range("I3").Value = 2400
$ws.Range("J3").Value = 2221.25
$ws.Range("K3").Value = 2400
$ws.Range("L3").Value = 2221.25
$ws.Range("M3").Value = -2286
$ws.Range("N3").Value = -2449.25
$ws.Range("H94").Value = 1578.0588
$ws.Range("I94").Value = 1388.4667
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 1388.4667
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -937.4666999999999
$ws.Range("N94").Value = -3902

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8448.125
$ws.Range("I31").Value = 1620
$ws.Range("J31").Value = 14472.941
$ws.Range("K31").Value = 1620
$ws.Range("L31").Value = 14472.941
$ws.Range("M31").Value = -1325
$ws.Range("N31").Value = -15062.941
$ws.Range("H34").Value = 8448.125
$ws.Range("I34").Value = 1620
$ws.Range("J34").Value = 14472.941
$ws.Range("K34").Value = 1620
$ws.Range("L34").Value = 14472.941
$ws.Range("M34").Value = -1418
$ws.Range("N34").Value = -14876.941
$ws.Range("H58").Value = 1588.4286
$ws.Range("I58").Value = 1154.9445
$ws.Range("J58").Value = 2368.7
$ws.Range("K58").Value = 1154.9445
$ws.Range("L58").Value = 2368.7
$ws.Range("M58").Value = -951.9445000000001
$ws.Range("N58").Value = -2774.7
$ws.Range("H132").Value = 3539.25
$ws.Range("I132").Value = 2807.375
$ws.Range("K132").Value = 8422.125
$ws.Range("M132").Value = -5892.125
$ws.Range("H134").Value = 1347.6945
$ws.Range("I134").Value = 1396.6552
$ws.Range("J134").Value = 1144.8572
$ws.Range("K134").Value = 4189.9656
$ws.Range("L134").Value = 3434.5716
$ws.Range("M134").Value = -1654.9656
$ws.Range("N134").Value = -8504.571599999999
$ws.Range("H136").Value = 1588.4286
$ws.Range("I136").Value = 1154.9445
$ws.Range("J136").Value = 2368.7
$ws.Range("K136").Value = 3464.8335
$ws.Range("L136").Value = 7106.099999999999
$ws.Range("M136").Value = -914.8335000000002
$ws.Range("N136").Value = -12206.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 316956.38
$ws.Range("I5").Value = 595.1429000000001
$ws.Range("K5").Value = 1785.4287
$ws.Range("M5").Value = -1673.4287
$ws.Range("H132").Value = 2146.606
$ws.Range("I132").Value = 1619.75
$ws.Range("K132").Value = 14577.75
$ws.Range("M132").Value = -12047.75
$ws.Range("H135").Value = 316956.38
$ws.Range("I135").Value = 595.1429000000001
$ws.Range("K135").Value = 5356.2861
$ws.Range("M135").Value = -2821.2861
$ws.Range("H137").Value = 55577530
$ws.Range("J137").Value = 166713310
$ws.Range("L137").Value = 500139930
$ws.Range("N137").Value = -500150130

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 200001220
$ws.Range("I113").Value = 1000000000
$ws.Range("J113").Value = 1525
$ws.Range("K113").Value = 1000000000
$ws.Range("L113").Value = 1525
$ws.Range("M113").Value = -999997830
$ws.Range("N113").Value = -5865

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H93").Value = 35730650
$ws.Range("I93").Value = 26849.75
$ws.Range("J93").Value = 83335710
$ws.Range("K93").Value = 26849.75
$ws.Range("L93").Value = 83335710
$ws.Range("M93").Value = -25601.75
$ws.Range("N93").Value = -83338206

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4557.7144
$ws.Range("I122").Value = 1725
$ws.Range("K122").Value = 5175
$ws.Range("M122").Value = -2725
